$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B header
$ws.Range("B1").Value = "Group"

# Column B data for rows 2..115 ("Group 1".."Group 6" assignment per person)
$groups = @("Group 5","Group 2","Group 3","Group 6","Group 5","Group 3","Group 1","Group 3","Group 3","Group 1","Group 4","Group 4","Group 4","Group 3","Group 4","Group 2","Group 4","Group 2","Group 1","Group 6","Group 1","Group 1","Group 4","Group 4","Group 5","Group 1","Group 5","Group 2","Group 3","Group 3","Group 3","Group 5","Group 5","Group 5","Group 1","Group 6","Group 5","Group 4","Group 1","Group 1","Group 3","Group 5","Group 1","Group 5","Group 6","Group 4","Group 1","Group 4","Group 6","Group 2","Group 1","Group 4","Group 1","Group 5","Group 4","Group 4","Group 3","Group 2","Group 4","Group 2","Group 3","Group 2","Group 4","Group 5","Group 5","Group 2","Group 5","Group 6","Group 1","Group 3","Group 3","Group 5","Group 2","Group 1","Group 2","Group 6","Group 2","Group 1","Group 3","Group 3","Group 2","Group 3","Group 6","Group 1","Group 6","Group 2","Group 6","Group 4","Group 6","Group 3","Group 4","Group 2","Group 5","Group 6","Group 4","Group 1","Group 6","Group 3","Group 6","Group 3","Group 4","Group 6","Group 5","Group 1","Group 1","Group 4","Group 6","Group 6","Group 2","Group 3","Group 2","Group 2","Group 5","Group 6")

for ($i = 0; $i -lt $groups.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $groups[$i]
}

# Widen column B to fit the new "Group N" labels
$ws.Columns("B").ColumnWidth = 16.28515625

# Update the active selection to match the saved view
$ws.Range("I25").Select()
